# "case with 380 kV done" -- refresh loading_percent results for Case_5_22.
# Rewrites the B2:N25 results block (columns B,C,D,E,F,H,J,L,M,N; G/I/K/O stay 0)
# with the recomputed loading percentages for rows 2-25 (A=0..23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'double[,]' 24,13
$data[0,0] = 29.22024917866318
$data[0,1] = 11.00672983195734
$data[0,2] = 3.364191789476627
$data[0,3] = 9.278120118370197
$data[0,4] = 57.5062078401759
$data[0,5] = 0
$data[0,6] = 7.344005520526261
$data[0,7] = 0
$data[0,8] = 9.681245611534154
$data[0,9] = 0
$data[0,10] = 11.66435812607003
$data[0,11] = 22.33059044669318
$data[0,12] = 23.80895137667823
$data[1,0] = 28.98779264211914
$data[1,1] = 10.67435376955665
$data[1,2] = 3.30457038433875
$data[1,3] = 9.263235039810688
$data[1,4] = 57.42628700465424
$data[1,5] = 0
$data[1,6] = 7.344005520526261
$data[1,7] = 0
$data[1,8] = 9.69077701673036
$data[1,9] = 0
$data[1,10] = 11.68308512607789
$data[1,11] = 22.30601294436337
$data[1,12] = 23.8646259361642
$data[2,0] = 28.85223907291062
$data[2,1] = 10.46883087925419
$data[2,2] = 3.266873852803161
$data[2,3] = 9.253889449250906
$data[2,4] = 57.38978736585815
$data[2,5] = 0
$data[2,6] = 7.344005520526261
$data[2,7] = 0
$data[2,8] = 9.696921398371403
$data[2,9] = 0
$data[2,10] = 11.69613837047859
$data[2,11] = 22.29553293496961
$data[2,12] = 23.90079462990975
$data[3,0] = 28.79885519298293
$data[3,1] = 10.38486423289621
$data[3,2] = 3.251244784431611
$data[3,3] = 9.250028918751992
$data[3,4] = 57.37807621702721
$data[3,5] = 0
$data[3,6] = 7.344005520526261
$data[3,7] = 0
$data[3,8] = 9.699498972866806
$data[3,9] = 0
$data[3,10] = 11.70184883001624
$data[3,11] = 22.29242433019914
$data[3,12] = 23.9160327916494
$data[4,0] = 28.79010431185187
$data[4,1] = 10.37091317956722
$data[4,2] = 3.248633628645783
$data[4,3] = 9.249384723398993
$data[4,4] = 57.37632260536536
$data[4,5] = 0
$data[4,6] = 7.344005520526261
$data[4,7] = 0
$data[4,8] = 9.699931435075955
$data[4,9] = 0
$data[4,10] = 11.7028206774359
$data[4,11] = 22.291978399294
$data[4,12] = 23.91859322781201
$data[5,0] = 28.85151154278624
$data[5,1] = 10.46769914173189
$data[5,2] = 3.26666414788081
$data[5,3] = 9.253837596224011
$data[5,4] = 57.38961661962462
$data[5,5] = 0
$data[5,6] = 7.344005520526261
$data[5,7] = 0
$data[5,8] = 9.696955861751029
$data[5,9] = 0
$data[5,10] = 11.69621379973292
$data[5,11] = 22.29548630309991
$data[5,12] = 23.90099811602124
$data[6,0] = 29.13864049072113
$data[6,1] = 10.89251240622492
$data[6,2] = 3.343864014672215
$data[6,3] = 9.27303033332274
$data[6,4] = 57.47604333462403
$data[6,5] = 0
$data[6,6] = 7.344005520526261
$data[6,7] = 0
$data[6,8] = 9.684471572823725
$data[6,9] = 0
$data[6,10] = 11.67049262073112
$data[6,11] = 22.32116101448043
$data[6,12] = 23.8277362199942
$data[7,0] = 29.75623174826032
$data[7,1] = 11.70804933666533
$data[7,2] = 3.486352041423892
$data[7,3] = 9.309046628936317
$data[7,4] = 57.7451227482041
$data[7,5] = 0
$data[7,6] = 7.344005520526261
$data[7,7] = 0
$data[7,8] = 9.662295882805234
$data[7,9] = 0
$data[7,10] = 11.63238140752491
$data[7,11] = 22.40794414869566
$data[7,12] = 23.69980657159582
$data[8,0] = 30.23979210309325
$data[8,1] = 12.28909187273435
$data[8,2] = 3.585286108132635
$data[8,3] = 9.334541539041759
$data[8,4] = 58.00318329665998
$data[8,5] = 0
$data[8,6] = 7.344005520526261
$data[8,7] = 0
$data[8,8] = 9.64739309749999
$data[8,9] = 0
$data[8,10] = 11.61188460691347
$data[8,11] = 22.49367241881287
$data[8,12] = 23.61539880506943
$data[9,0] = 30.46547013018498
$data[9,1] = 12.54811432529454
$data[9,2] = 3.628982380617649
$data[9,3] = 9.345934723957274
$data[9,4] = 58.13356501851267
$data[9,5] = 0
$data[9,6] = 7.344005520526261
$data[9,7] = 0
$data[9,8] = 9.640911753624176
$data[9,9] = 0
$data[9,10] = 11.60418633049801
$data[9,11] = 22.53737493345056
$data[9,12] = 23.5790782899205
$data[10,0] = 30.55168059710559
$data[10,1] = 12.64533743405289
$data[10,2] = 3.645335868425309
$data[10,3] = 9.350219944383184
$data[10,4] = 58.1847900095064
$data[10,5] = 0
$data[10,6] = 7.344005520526261
$data[10,7] = 0
$data[10,8] = 9.638500029779944
$data[10,9] = 0
$data[10,10] = 11.60150466037359
$data[10,11] = 22.5545933235548
$data[10,12] = 23.56562331700717
$data[11,0] = 30.53308133786454
$data[11,1] = 12.62443859125053
$data[11,2] = 3.641822528191019
$data[11,3] = 9.349298336576396
$data[11,4] = 58.17367568376214
$data[11,5] = 0
$data[11,6] = 7.344005520526261
$data[11,7] = 0
$data[11,8] = 9.639017545800503
$data[11,9] = 0
$data[11,10] = 11.60207182604467
$data[11,11] = 22.55085539640985
$data[11,12] = 23.56850779343951
$data[12,0] = 30.47254806408126
$data[12,1] = 12.55613077130129
$data[12,2] = 3.630331694558228
$data[12,3] = 9.346287853819515
$data[12,4] = 58.13774230465755
$data[12,5] = 0
$data[12,6] = 7.344005520526261
$data[12,7] = 0
$data[12,8] = 9.640712486804958
$data[12,9] = 0
$data[12,10] = 11.60396102998656
$data[12,11] = 22.53877813346196
$data[12,12] = 23.57796535149244
$data[13,0] = 30.43556542960549
$data[13,1] = 12.51417502408237
$data[13,2] = 3.623267903843463
$data[13,3] = 9.344440060920112
$data[13,4] = 58.11597279022862
$data[13,5] = 0
$data[13,6] = 7.344005520526261
$data[13,7] = 0
$data[13,8] = 9.641756230657188
$data[13,9] = 0
$data[13,10] = 11.60514862131341
$data[13,11] = 22.53146738133452
$data[13,12] = 23.58379729870317
$data[14,0] = 30.22515241031899
$data[14,1] = 12.27204866277466
$data[14,2] = 3.58240366644949
$data[14,3] = 9.333792875516101
$data[14,4] = 57.9949226413108
$data[14,5] = 0
$data[14,6] = 7.344005520526261
$data[14,7] = 0
$data[14,8] = 9.64782264532338
$data[14,9] = 0
$data[14,10] = 11.61242039898441
$data[14,11] = 22.49091041469266
$data[14,12] = 23.61781424262071
$data[15,0] = 30.09748298281171
$data[15,1] = 12.12207998800824
$data[15,2] = 3.556995720616771
$data[15,3] = 9.327208979235749
$data[15,4] = 57.92397876272241
$data[15,5] = 0
$data[15,6] = 7.344005520526261
$data[15,7] = 0
$data[15,8] = 9.651620355230165
$data[15,9] = 0
$data[15,10] = 11.61729761488003
$data[15,11] = 22.46722997937755
$data[15,12] = 23.63921460068998
$data[16,0] = 30.02459202233177
$data[16,1] = 12.03532883014726
$data[16,2] = 3.542258783062527
$data[16,3] = 9.323402823371897
$data[16,4] = 57.88439736140336
$data[16,5] = 0
$data[16,6] = 7.344005520526261
$data[16,7] = 0
$data[16,8] = 9.653832761289614
$data[16,9] = 0
$data[16,10] = 11.62025591009101
$data[16,11] = 22.45405289460652
$data[16,12] = 23.65171900709731
$data[17,0] = 30.00000741902499
$data[17,1] = 12.00587493721658
$data[17,2] = 3.537248150674005
$data[17,3] = 9.322110805521186
$data[17,4] = 57.87120637239654
$data[17,5] = 0
$data[17,6] = 7.344005520526261
$data[17,7] = 0
$data[17,8] = 9.654586671052586
$data[17,9] = 0
$data[17,10] = 11.62128383271306
$data[17,11] = 22.44966770212439
$data[17,12] = 23.65598635520175
$data[18,0] = 30.11101810472458
$data[18,1] = 12.13809621457039
$data[18,2] = 3.559713203364801
$data[18,3] = 9.327911838821041
$data[18,4] = 57.93140432200774
$data[18,5] = 0
$data[18,6] = 7.344005520526261
$data[18,7] = 0
$data[18,8] = 9.65121317972447
$data[18,9] = 0
$data[18,10] = 11.61676258886745
$data[18,11] = 22.46970497511639
$data[18,12] = 23.63691626297409
$data[19,0] = 30.49030832592291
$data[19,1] = 12.57621861565175
$data[19,2] = 3.633712116859682
$data[19,3] = 9.347172893922824
$data[19,4] = 58.14824667478292
$data[19,5] = 0
$data[19,6] = 7.344005520526261
$data[19,7] = 0
$data[19,8] = 9.640213486689371
$data[19,9] = 0
$data[19,10] = 11.60339979036247
$data[19,11] = 22.54230741830751
$data[19,12] = 23.57517932626489
$data[20,0] = 30.74254161088959
$data[20,1] = 12.85748626815439
$data[20,2] = 3.680945628575554
$data[20,3] = 9.359591438815478
$data[20,4] = 58.30075320187787
$data[20,5] = 0
$data[20,6] = 7.344005520526261
$data[20,7] = 0
$data[20,8] = 9.633272868601731
$data[20,9] = 0
$data[20,10] = 11.59602728615274
$data[20,11] = 22.59365399215798
$data[20,12] = 23.53657249984479
$data[21,0] = 30.60754548081957
$data[21,1] = 12.70786331414567
$data[21,2] = 3.655841108120077
$data[21,3] = 9.352978857944715
$data[21,4] = 58.21837628340446
$data[21,5] = 0
$data[21,6] = 7.344005520526261
$data[21,7] = 0
$data[21,8] = 9.636954562669553
$data[21,9] = 0
$data[21,10] = 11.59983771714622
$data[21,11] = 22.56589541208287
$data[21,12] = 23.55701825493155
$data[22,0] = 30.10489728890202
$data[22,1] = 12.13085692832549
$data[22,2] = 3.558485032611604
$data[22,3] = 9.327594141391932
$data[22,4] = 57.92804347164299
$data[22,5] = 0
$data[22,6] = 7.344005520526261
$data[22,7] = 0
$data[22,8] = 9.651397173262541
$data[22,9] = 0
$data[22,10] = 11.61700399343748
$data[22,11] = 22.46858466749095
$data[22,12] = 23.63795471522323
$data[23,0] = 29.58366235797743
$data[23,1] = 11.49010924548661
$data[23,2] = 3.448792992601394
$data[23,3] = 9.299474512422909
$data[23,4] = 57.66169323634143
$data[23,5] = 0
$data[23,6] = 7.344005520526261
$data[23,7] = 0
$data[23,8] = 9.668049792046601
$data[23,9] = 0
$data[23,10] = 11.64137281694567
$data[23,11] = 22.38058834222004
$data[23,12] = 23.73273149448173

$ws.Range("B2:N25").Value = $data
